# Adjusted excel + jupyter for start-up/shut-down limits
# Adds 4 new columns (start_up_Output1, start_up_Output2, shut_down_Output1,
# shut_down_Output2) to the "Units" table, positioned right after
# "ramp_down_Output2" and before "Relation_In_In", and populates the new
# start-up / shut-down limit values for the Electrolyzer row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Units")

# Remove the existing table definition but keep its data in place so the
# underlying cells can be shifted with a plain column insert.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Unlist()

# Insert 4 blank columns right before "Relation_In_In" (was column S),
# shifting Relation_In_In .. vom_cost_Output2 four columns to the right
# (S:V -> W:Z, etc).
$ws.Columns("S:V").Insert()

# New column headers.
$ws.Range("S1").Value2 = "start_up_Output1"
$ws.Range("T1").Value2 = "start_up_Output2"
$ws.Range("U1").Value2 = "shut_down_Output1"
$ws.Range("V1").Value2 = "shut_down_Output2"

# New data values for the Electrolyzer row (row 2).
$ws.Range("S2").Value2 = 0.5
$ws.Range("U2").Value2 = 0.8

# Re-create the table over the full, now-wider range and restore its name.
$newtbl = $ws.ListObjects.Add(1, $ws.Range("A1:AG6"), 0, 1)
$newtbl.Name = "Table1"
